# Apply "Added more header related tests" changes to the API_Template sheet
# of the QA API doc workbook:
#  - Row 11 ("Create task") becomes "Create a task" with a new endpoint and a
#    bigger JSON payload template (row grows taller to fit the payload).
#  - A brand-new row 14 is inserted for "Get all tasks" (GET on the same
#    /tasks endpoint).
#  - The sheet's last active-cell selection moves to B17.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("API_Template")

# --- Row 11: "Create task" -> "Create a task" -------------------------------
$ws.Range("A11").Value2 = "Create a task"
$ws.Range("B11").Value2 = "https://intelliapi-mockserver.herokuapp.com/tasks"

$newPayload = @"
{
    "status": [
        "#status"
    ],
    "name": "#name",
    "category": "#category",
    "isDeleted": #isDeleted,
    "__v": #version
}
"@
$ws.Range("E11").Value2 = $newPayload

# Row grows from 45pt to 135pt to fit the larger JSON template.
$ws.Rows.Item(11).RowHeight = 135

# Point the existing B11 hyperlink at the new endpoint (best effort - the
# original hyperlink loaded from the workbook can't be removed by this host,
# so re-pointing/adding is the supported path).
$ws.Range("B11").Hyperlinks.Item(1).Address = "https://intelliapi-mockserver.herokuapp.com/tasks"

# --- New row 14: "Get all tasks" --------------------------------------------
# Copy row 12's formatting (same visual style as the other "GET" rows) into
# the new row, then fill in the new content.
$ws.Range("A12:E12").Copy()
$ws.Range("A14:E14").PasteSpecial(-4122)
$ws.Rows.Item(14).RowHeight = 15

$ws.Range("A14").Value2 = "Get all tasks"
$ws.Range("B14").Value2 = "https://intelliapi-mockserver.herokuapp.com/tasks"
$ws.Range("C14").Value2 = "GET"

$ws.Hyperlinks.Add($ws.Range("B14"), "https://intelliapi-mockserver.herokuapp.com/tasks")

# --- Sheet view: last selected cell moves to B17 ----------------------------
$ws.Range("B17").Select()
